# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> clrScheme "Integral"      (used by the slide master)
#   ppt/theme/theme2.xml -> clrScheme "Office"        (used by the notes master)
#
# The target edit swaps the two palettes: the slide master's theme becomes the
# stock "Office Theme" palette (what used to live in theme2.xml) while the
# notes-master theme keeps/reverts to the "Integral" palette.
#
# The PowerPoint object model only exposes a single addressable Theme (the one
# driving SlideMaster/Slides/Layouts), so we rewrite its twelve theme colors,
# in clrScheme order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), to the
# Office Theme RGB values. ThemeColorScheme.Colors(i).RGB uses the VBA/OLE
# &H00BBGGRR encoding, i.e. the bytes of the RRGGBB hex code reversed.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# index -> RRGGBB for the stock "Office" color scheme
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $oleColor = ($bb * 65536) + ($gg * 256) + $rr
    $colorScheme.Colors($i).RGB = $oleColor
}
